$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A43").Value = "Backend Engineer (GoLang & TypeScript)"
$ws.Range("B43").Value = "https://www.dice.com/job-detail/0ec930c6-7d7c-452a-87bf-36d2d64b9b4c"
$ws.Range("C43").Value = "Remote or New Jersey"
$ws.Range("D43").Value = "Full-time, Contract"
$ws.Range("E43").Value = "Depends on Experience"
$ws.Range("F43").Value = "Radyant Inc."

$ws.Range("A44").Value = "Sr. Golang Developer- (Cloud Software Engineer)"
$ws.Range("B44").Value = "https://www.dice.com/job-detail/8b7471d3-f952-44db-bd51-b67cc17cdd03"
$ws.Range("C44").Value = "Hybrid in Plano, Texas"
$ws.Range("D44").Value = "Contract"
$ws.Range("E44").Value = "100 - 110"
$ws.Range("F44").Value = "MSYS Inc."
